$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Cells.Item(2, 1).Value = 580644304
$ws.Cells.Item(2, 2).Value = 3

# Remove rows 3 through 8 (old data rows no longer present)
$ws.Range("A3:B8").ClearContents()
